$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume columns so values are not
# auto-converted to numbers by the Excel value-assignment heuristics.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.687.57"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "1.891.38"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "245.23"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4921"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.2962"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("D9").Value = "0.06793"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("D10").Value = "1.886.82"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").Value = "17.30"
$ws.Range("E11").Value = "  +4.04%  "
$ws.Range("D12").Value = "0.07229"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "91.40"
$ws.Range("E13").Value = "  +6.09%  "
$ws.Range("D14").Value = "0.6787"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "5.050"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("D16").Value = "30.636.54"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").Value = "0.000008000"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "0.9999"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "2.130.04"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "4.821"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "193.37"
$ws.Range("E23").Value = "  +35.23%  "
$ws.Range("D24").Value = "6.098"
$ws.Range("E24").Value = "  +4.24%  "
$ws.Range("D25").Value = "9.365"
$ws.Range("E25").Value = "  +3.11%  "
$ws.Range("D26").Value = "155.07"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").Value = "19.27"
$ws.Range("E27").Value = "  +13.75%  "
$ws.Range("D28").Value = "1.908"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "1.404"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").Value = "4.346"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").Value = "0.09099"
$ws.Range("E31").Value = "  +4.03%  "
$ws.Range("D32").Value = "4.018"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "0.05207"
$ws.Range("E33").Value = "  +4.02%  "
$ws.Range("D34").Value = "0.7634"
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("D35").Value = "1.112"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "2.778"
$ws.Range("E36").Value = "  +4.19%  "
$ws.Range("D37").Value = "0.01847"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("D38").Value = "2.680"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "2.150"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "0.9356"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "0.4433"
$ws.Range("E41").Value = "  +5.12%  "
$ws.Range("D42").Value = "105.50"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D45").Value = "7.614"
$ws.Range("E45").Value = "  +3.71%  "
$ws.Range("D46").Value = "0.1347"
$ws.Range("E46").Value = "  +6.17%  "
$ws.Range("D47").Value = "0.05864"
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("D48").Value = "8.717"
$ws.Range("E48").Value = "  +5.43%  "
$ws.Range("E49").Value = "  +6.61%  "
$ws.Range("D50").Value = "0.3937"
$ws.Range("E50").Value = "  +4.65%  "
$ws.Range("D51").Value = "33.64"
$ws.Range("E51").Value = "  +2.87%  "

# Row 43/44 swap: PaxDollar <-> FraxShare with updated data
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.770"
$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  +0.09%  "

# Restore default (General) styling so no stray number-format is left behind
$ws.Range("B2:E51").Style = "Normal"

Write-Host "done"
